# LoginData workbook update: added Extent Report support, added Listener classes.
# Functional data change: header renamed Mail_Id -> Username, email/password
# test data replaced with a shorter 5-row "testmeN@gmail.com" / passNNN set,
# mail addresses turned into mailto: hyperlinks, and the now-unused rows
# 7-10 cleared out (while keeping their existing formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# --- Header row ---
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

# --- Data rows 2-6: new email / password pairs ---
$ws.Range("A2").Value = "testme1@gmail.com"
$ws.Range("B2").Value = "pass123"

$ws.Range("A3").Value = "testme2@gmail.com"
$ws.Range("B3").Value = "pass124"

$ws.Range("A4").Value = "testme3@gmail.com"
$ws.Range("B4").Value = "pass125"

$ws.Range("A5").Value = "testme4@gmail.com"
$ws.Range("B5").Value = "pass126"

$ws.Range("A6").Value = "testme5@gmail.com"
$ws.Range("B6").Value = "pass127"

# --- Clear the now-unused trailing rows, keep their formatting ---
$ws.Range("A7:B10").ClearContents()

# --- Turn the e-mail cells into mailto: hyperlinks ---
# (Add() order reproduces the relationship-id ordering recorded in the
# saved file: A3, A2, A5, A4, A6.)
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:testme2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:testme1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:testme4@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:testme3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:testme5@gmail.com")

# Hyperlinks.Add() re-styles the target cells with a freshly minted xf;
# restore the workbook's existing built-in "Hyperlink" style so the cells
# keep using the same style index as before.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("A5").Style = "Hyperlink"
$ws.Range("A6").Style = "Hyperlink"
